# --- Part 1: update "总计" (summary) sheet - insert 2022-Q3 row ---
$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item(1)
$wsSummary.Rows.Item(2).Insert()

$wsSummary = $wb.Worksheets.Item(1)
$wsSummary.Range("A3").Copy()
$wsSummary.Range("A2").PasteSpecial(-4122)
$wsSummary.Range("B2:D2").ClearFormats()

$wsSummary = $wb.Worksheets.Item(1)
$wsSummary.Range("A2").Value = 7
$wsSummary.Range("B2").Value = "2022-Q3"
$wsSummary.Range("C2").Value = 18
$wsSummary.Range("D2").Value = 2.93

# --- Part 2: create the new "2022-Q3" fund-holdings sheet before "2022-Q2" ---
$target = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "2022-Q3"

# Copy header-row style (bold/border, s=2) from the "2022-Q2" sheet's header row.
$srcHeader = $wb.Worksheets.Item("2022-Q2")
$srcHeader.Range("B1:H1").Copy()
$dstHeader = $wb.Worksheets.Item("2022-Q3")
$dstHeader.Range("B1:H1").PasteSpecial(-4122)

# Copy the column-A index style (s=2) down for all 18 data rows.
$srcIdx = $wb.Worksheets.Item("2022-Q2")
$srcIdx.Range("A2:H2").Copy()
$dstIdx = $wb.Worksheets.Item("2022-Q3")
$dstIdx.Range("A2:H19").PasteSpecial(-4122)

# --- Header row text for the new sheet ---
$ws3 = $wb.Worksheets.Item("2022-Q3")
$ws3.Range("B1").Value = "基金代码"
$ws3.Range("C1").Value = "基金名称"
$ws3.Range("D1").Value = "基金规模"
$ws3.Range("E1").Value = "股票总仓位"
$ws3.Range("F1").Value = "仓位占比"
$ws3.Range("G1").Value = "持有市值(亿元)"
$ws3.Range("H1").Value = "仓位排名"

# --- Data rows for the new sheet ---
$ws3 = $wb.Worksheets.Item("2022-Q3")
$ws3.Range("B2:G19").NumberFormat = "@"
$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "002363"
$ws3.Range("C2").Value = "华安安康灵活配置混合A"
$ws3.Range("D2").Value = "107.76"
$ws3.Range("E2").Value = "21.95"
$ws3.Range("F2").Value = "1.40"
$ws3.Range("G2").Value = "1.5086"
$ws3.Range("H2").Value = 3
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "002364"
$ws3.Range("C3").Value = "华安安康灵活配置混合C"
$ws3.Range("D3").Value = "27.33"
$ws3.Range("E3").Value = "21.95"
$ws3.Range("F3").Value = "1.40"
$ws3.Range("G3").Value = "0.3826"
$ws3.Range("H3").Value = 3
$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "002144"
$ws3.Range("C4").Value = "华安新优选灵活配置混合C"
$ws3.Range("D4").Value = "19.79"
$ws3.Range("E4").Value = "22.13"
$ws3.Range("F4").Value = "1.23"
$ws3.Range("G4").Value = "0.2434"
$ws3.Range("H4").Value = 4
$ws3.Range("A5").Value = 3
$ws3.Range("B5").Value = "005695"
$ws3.Range("C5").Value = "华安睿明两年定期开放灵活配置混合A"
$ws3.Range("D5").Value = "4.27"
$ws3.Range("E5").Value = "93.55"
$ws3.Range("F5").Value = "4.23"
$ws3.Range("G5").Value = "0.1806"
$ws3.Range("H5").Value = 5
$ws3.Range("A6").Value = 4
$ws3.Range("B6").Value = "001312"
$ws3.Range("C6").Value = "华安新优选灵活配置混合A"
$ws3.Range("D6").Value = "9.83"
$ws3.Range("E6").Value = "22.13"
$ws3.Range("F6").Value = "1.23"
$ws3.Range("G6").Value = "0.1209"
$ws3.Range("H6").Value = 4
$ws3.Range("A7").Value = 5
$ws3.Range("B7").Value = "040020"
$ws3.Range("C7").Value = "华安升级主题混合A"
$ws3.Range("D7").Value = "4.40"
$ws3.Range("E7").Value = "85.97"
$ws3.Range("F7").Value = "2.70"
$ws3.Range("G7").Value = "0.1188"
$ws3.Range("H7").Value = 9
$ws3.Range("A8").Value = 6
$ws3.Range("B8").Value = "160425"
$ws3.Range("C8").Value = "华安创业板两年定期开放混合"
$ws3.Range("D8").Value = "1.80"
$ws3.Range("E8").Value = "93.72"
$ws3.Range("F8").Value = "4.84"
$ws3.Range("G8").Value = "0.0871"
$ws3.Range("H8").Value = 9
$ws3.Range("A9").Value = 7
$ws3.Range("B9").Value = "001028"
$ws3.Range("C9").Value = "华安物联网主题股票A"
$ws3.Range("D9").Value = "2.96"
$ws3.Range("E9").Value = "94.04"
$ws3.Range("F9").Value = "2.89"
$ws3.Range("G9").Value = "0.0855"
$ws3.Range("H9").Value = 2
$ws3.Range("A10").Value = 8
$ws3.Range("B10").Value = "011390"
$ws3.Range("C10").Value = "华安添祥6个月持有期混合A"
$ws3.Range("D10").Value = "6.67"
$ws3.Range("E10").Value = "33.77"
$ws3.Range("F10").Value = "1.19"
$ws3.Range("G10").Value = "0.0794"
$ws3.Range("H10").Value = 7
$ws3.Range("A11").Value = 9
$ws3.Range("B11").Value = "008290"
$ws3.Range("C11").Value = "华安现代生活混合"
$ws3.Range("D11").Value = "2.68"
$ws3.Range("E11").Value = "90.97"
$ws3.Range("F11").Value = "2.73"
$ws3.Range("G11").Value = "0.0732"
$ws3.Range("H11").Value = 9
$ws3.Range("A12").Value = 10
$ws3.Range("B12").Value = "000549"
$ws3.Range("C12").Value = "华安大国新经济股票A"
$ws3.Range("D12").Value = "0.72"
$ws3.Range("E12").Value = "93.83"
$ws3.Range("F12").Value = "3.82"
$ws3.Range("G12").Value = "0.0275"
$ws3.Range("H12").Value = 1
$ws3.Range("A13").Value = 11
$ws3.Range("B13").Value = "009409"
$ws3.Range("C13").Value = "华安添福18个月持有期混合A"
$ws3.Range("D13").Value = "0.72"
$ws3.Range("E13").Value = "21.56"
$ws3.Range("F13").Value = "1.91"
$ws3.Range("G13").Value = "0.0138"
$ws3.Range("H13").Value = 3
$ws3.Range("A14").Value = 12
$ws3.Range("B14").Value = "005696"
$ws3.Range("C14").Value = "华安睿明两年定期开放灵活配置混合C"
$ws3.Range("D14").Value = "0.07"
$ws3.Range("E14").Value = "93.55"
$ws3.Range("F14").Value = "4.23"
$ws3.Range("G14").Value = "0.0030"
$ws3.Range("H14").Value = 5
$ws3.Range("A15").Value = 13
$ws3.Range("B15").Value = "009410"
$ws3.Range("C15").Value = "华安添福18个月持有期混合C"
$ws3.Range("D15").Value = "0.09"
$ws3.Range("E15").Value = "21.56"
$ws3.Range("F15").Value = "1.91"
$ws3.Range("G15").Value = "0.0017"
$ws3.Range("H15").Value = 3
$ws3.Range("A16").Value = 14
$ws3.Range("B16").Value = "014976"
$ws3.Range("C16").Value = "华安升级主题混合C"
$ws3.Range("D16").Value = "0.00"
$ws3.Range("E16").Value = "85.97"
$ws3.Range("F16").Value = "2.70"
$ws3.Range("G16").NumberFormat = "General"
$ws3.Range("G16").Value = 0
$ws3.Range("H16").Value = 9
$ws3.Range("A17").Value = 15
$ws3.Range("B17").Value = "016291"
$ws3.Range("C17").Value = "华安大国新经济股票C"
$ws3.Range("D17").Value = "0.00"
$ws3.Range("E17").Value = "93.83"
$ws3.Range("F17").Value = "3.82"
$ws3.Range("G17").NumberFormat = "General"
$ws3.Range("G17").Value = 0
$ws3.Range("H17").Value = 1
$ws3.Range("A18").Value = 16
$ws3.Range("B18").Value = "016292"
$ws3.Range("C18").Value = "华安物联网主题股票C"
$ws3.Range("D18").Value = "0.00"
$ws3.Range("E18").Value = "94.04"
$ws3.Range("F18").Value = "2.89"
$ws3.Range("G18").NumberFormat = "General"
$ws3.Range("G18").Value = 0
$ws3.Range("H18").Value = 2
$ws3.Range("A19").Value = 17
$ws3.Range("B19").Value = "016181"
$ws3.Range("C19").Value = "华安添祥6个月持有期混合C"
$ws3.Range("D19").Value = "0.00"
$ws3.Range("E19").Value = "33.77"
$ws3.Range("F19").Value = "1.19"
$ws3.Range("G19").NumberFormat = "General"
$ws3.Range("G19").Value = 0
$ws3.Range("H19").Value = 7
